# Updates cryptos list data (price + 1h volume change columns) to match the
# latest scrape, and fixes the Maker / HuobiToken row ordering (rows 35-36).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.654.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.98%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.632.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.39%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.43%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.498'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.75%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.43%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.93%  '

# Row 9
$ws.Range("E9").Value = '  +0.75%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.34%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.52%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.859.99'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.43%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.610.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.94%  '

# Row 14
$ws.Range("E14").Value = '  +1.33%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.07%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.646.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.01%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.56'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.42%  '

# Row 18
$ws.Range("E18").Value = '  +1.85%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.01%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.45%  '

# Row 21
$ws.Range("E21").Value = '  +0.80%  '

# Row 22
$ws.Range("E22").Value = '  +1.56%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.45%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.18%  '

# Row 26
$ws.Range("E26").Value = '  +0.33%  '

# Row 27
$ws.Range("E27").Value = '  +1.15%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.83%  '

# Row 29
$ws.Range("E29").Value = '  +2.13%  '

# Row 31
$ws.Range("E31").Value = '  -0.18%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.56%  '

# Row 34
$ws.Range("E34").Value = '  +0.33%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.40'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.08%  '

# Row 36
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.224.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.59%  '

# Row 37
$ws.Range("E37").Value = '  +5.50%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.803'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.31%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.38%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.500'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.88%  '

# Row 41
$ws.Range("E41").Value = '  -1.45%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.797'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.75%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.92%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.770.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.41%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.00%  '

# Row 46
$ws.Range("E46").Value = '  +2.51%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.09%  '

# Row 48
$ws.Range("E48").Value = '  +0.80%  '

# Row 49
$ws.Range("E49").Value = '  +1.03%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.59%  '

# Row 51
$ws.Range("E51").Value = '  +0.17%  '
